# Update cryptos list with refreshed prices / volume(1h) percentages
# and restore the two swapped coin-name/link pairs (rows 14/15, 19/20, 21/22, 31/33)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "84.403.71"
$ws.Range("E2").Value = "  +5.74%  "
$ws.Range("D3").Value = "3.289.81"
$ws.Range("E3").Value = "  +2.33%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'217.90"
$ws.Range("E5").Value = "  +3.05%  "
$ws.Range("D6").Value = "'633.94"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D7").Value = "'0.322"
$ws.Range("E7").Value = "  +20.54%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.592"
$ws.Range("E9").Value = "  -2.38%  "
$ws.Range("D10").Value = "3.289.10"
$ws.Range("E10").Value = "  +2.48%  "
$ws.Range("D11").Value = "'0.593"
$ws.Range("E11").Value = "  -4.02%  "
$ws.Range("D12").Value = "'0.0000278"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.892.91"
$ws.Range("E14").Value = "  +2.41%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'34.01"
$ws.Range("E15").Value = "  +3.71%  "
$ws.Range("D16").Value = "'5.43"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "84.370.84"
$ws.Range("E17").Value = "  +5.98%  "
$ws.Range("D18").Value = "3.292.33"
$ws.Range("E18").Value = "  +3.37%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'14.56"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D20").Value = "'3.19"
$ws.Range("E20").Value = "  +5.15%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'9.18"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "'432.93"
$ws.Range("E22").Value = "  -3.23%  "
$ws.Range("D23").Value = "'5.20"
$ws.Range("E23").Value = "  -2.79%  "
$ws.Range("E24").Value = "  +5.43%  "
$ws.Range("D25").Value = "'5.45"
$ws.Range("E25").Value = "  +11.61%  "
$ws.Range("D26").Value = "'12.13"
$ws.Range("E26").Value = "  +10.54%  "
$ws.Range("D27").Value = "3.458.11"
$ws.Range("E27").Value = "  +2.73%  "
$ws.Range("D28").Value = "'78.05"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").Value = "'0.0000131"
$ws.Range("E29").Value = "  +3.36%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").Value = "'0.163"
$ws.Range("E31").Value = "  +32.13%  "
$ws.Range("D32").Value = "'597.26"
$ws.Range("E32").Value = "  +5.69%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'9.28"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("E35").Value = "  +2.31%  "
$ws.Range("D36").Value = "'0.153"
$ws.Range("E36").Value = "  -3.37%  "
$ws.Range("E37").Value = "  -1.18%  "
$ws.Range("D38").Value = "'23.20"
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("D39").Value = "'6.37"
$ws.Range("E39").Value = "  +9.34%  "
$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").Value = "'3.08"
$ws.Range("E42").Value = "  +11.87%  "
$ws.Range("D43").Value = "'2.05"
$ws.Range("E43").Value = "  +11.81%  "
$ws.Range("D44").Value = "'20.93"
$ws.Range("E44").Value = "  +3.06%  "
$ws.Range("D45").Value = "'158.99"
$ws.Range("E45").Value = "  -2.90%  "
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "'189.94"
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D48").Value = "'45.06"
$ws.Range("E48").Value = "  +4.46%  "
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").Value = "'0.783"
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("D51").Value = "'26.59"
$ws.Range("E51").Value = "  +2.26%  "
